$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.318.97"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.591.92"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'211.99"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "'0.502"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "'19.38"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "1.817.23"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.04"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.553.97"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "'64.52"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "26.334.17"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "'7.48"
$ws.Range("E19").Value = "  +3.29%  "
$ws.Range("D20").Value = "'211.80"
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "'4.29"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "'9.00"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = "  -1.75%  "
$ws.Range("D25").Value = "'144.88"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").Value = "'15.19"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "'0.0501"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").Value = "'2.97"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").Value = "1.338.21"
$ws.Range("E34").Value = "  +4.44%  "
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").Value = "'0.603"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "'1.48"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").Value = "'0.0166"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  -15.90%  "
$ws.Range("D40").Value = "'0.818"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").Value = "'5.77"
$ws.Range("E41").Value = "  +4.70%  "
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").Value = "'0.763"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "1.729.37"
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("D46").Value = "'61.88"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("E48").Value = "  +7.93%  "
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("D50").Value = "'0.0982"
$ws.Range("E50").Value = "  -2.62%  "
$ws.Range("E51").Value = "  -0.48%  "
